$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: "TV do banheiro" (TV in the bathroom) — correction noted in the
# commit message ("Correção da TV na planilha"). Columns A, B, D, E are
# populated; C and F are left empty, matching the other rows in the sheet.

# Columns A, B and D carry a default column style (see <cols> in the sheet
# XML). Typing straight into row 11 would make Excel inherit that column
# style for the new cells, but every other data row (2-10) in this sheet
# uses the workbook's default (unstyled) formatting. Copy the plain format
# from an already-unstyled cell in the same column first so the new row
# matches the rest of the table, then fill in the values.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A11").Value = "TV do banheiro"
$ws.Range("B11").Value = "televisao"
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = $false
